$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing row ---
$ws.Range("A2").Value = "044/FES VILLE /AV6"
$ws.Range("H2").Value = "--"
$ws.Range("J2").Value = "--"
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 600
$ws.Range("N2").Value = "--"
$ws.Range("O2").Value = 9400

# --- Row 3: was the blank totals row, now becomes a data row ---
$ws.Range("A3").Value = "044/FES VILLE /AV6"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "K5443645"
$ws.Range("D3").Value = "KHADIJA LALA"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 15
$ws.Range("H3").Value = 30000
$ws.Range("J3").Value = 4500
$ws.Range("O3").Value = 25500

# --- Row 4: brand-new data row ---
$ws.Range("A4").Value = "044/FES VILLE /AV6"
$ws.Range("B4").Value = "Direction régionale"
$ws.Range("C4").Value = "K5443645"
$ws.Range("D4").Value = "KHADIJA LALA"
$ws.Range("E4").Value = "non"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = "--"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "--"
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = 700
$ws.Range("N4").Value = "--"
$ws.Range("O4").Value = 4300

# --- Row 5: brand-new totals row (blank labels, totals in numeric columns) ---
$ws.Range("A5").Value = " "
$ws.Range("B5").Value = " "
$ws.Range("C5").Value = " "
$ws.Range("D5").Value = " "
$ws.Range("E5").Value = " "
$ws.Range("F5").Value = " "
$ws.Range("G5").Value = " "
$ws.Range("H5").Value = 30000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 4500
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 1300
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 39200
